# Wherewhenhownofoa_80_firstAttempt experiment
# - flips six of the Correct/incorrect input bits on row 6 (this drives every
#   downstream BKT formula in rows 7-11 via automatic recalculation)
# - logs two snapshots of the resulting Mastery sequence (rows 15 & 17) as
#   plain pasted values (their previous cell formatting is cleared, matching
#   a values-only paste)
# - leaves the selection on the just-pasted C17:I17 range

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- toggle the sequence of correct/incorrect answers (row 6) ---
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 1
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 1
$ws.Range("O6").Value = 0
$ws.Range("Q6").Value = 1

# --- log row 15: snapshot of a mastery-prediction run, pasted as values ---
$row15 = @{
    "C" = 0.1310124970296011
    "D" = 0.14668493374016453
    "E" = 0.65148849982329493
    "F" = 0.94959590801312577
    "G" = 0.99472299318753021
    "H" = 0.99946969415598019
    "I" = 0.99623161876913369
    "J" = 0.97383021775343137
    "K" = 0.84295215286624348
    "L" = 0.999730253533889
}
foreach ($col in $row15.Keys) {
    $cell = $ws.Range($col + "15")
    $cell.ClearFormats()
    $cell.Value = $row15[$col]
}

# --- log row 17: snapshot of the (now recalculated) Mastery row, pasted as values ---
$row17 = @{
    "C" = 0.1310124970296011
    "D" = 0.14668493374016453
    "E" = 0.14883697341285768
    "F" = 0.65504969461346618
    "G" = 0.95033723993703234
    "H" = 0.73913487150550528
    "I" = 0.96607620411026762
}
foreach ($col in $row17.Keys) {
    $cell = $ws.Range($col + "17")
    $cell.ClearFormats()
    $cell.Value = $row17[$col]
}

# --- leave the selection where the user finished: the pasted C17:I17 range ---
$ws.Range("C17:I17").Select()
